# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# to reflect the refreshed crypto data snapshot (GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'29.869.78"
$c.ClearFormats()
$ws.Range('E2').Value = '  +1.06%  '

$c = $ws.Range('D3')
$c.Value = "'1.622.41"
$c.ClearFormats()
$ws.Range('E3').Value = '  +1.19%  '

$c = $ws.Range('D4')
$c.Value = "'0.994"
$c.ClearFormats()
$ws.Range('E4').Value = '  -0.41%  '

$c = $ws.Range('D5')
$c.Value = "'213.70"
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.62%  '

$ws.Range('E6').Value = '  -0.35%  '

$c = $ws.Range('D7')
$c.Value = "'0.994"
$c.ClearFormats()
$ws.Range('E7').Value = '  -0.43%  '

$c = $ws.Range('D8')
$c.Value = "'29.82"
$c.ClearFormats()
$ws.Range('E8').Value = '  +11.38%  '

$ws.Range('E9').Value = '  +3.24%  '

$c = $ws.Range('D10')
$c.Value = "'0.0609"
$c.ClearFormats()
$ws.Range('E10').Value = '  +1.34%  '

$ws.Range('E11').Value = '  +0.56%  '

$c = $ws.Range('D12')
$c.Value = "'1.855.91"
$c.ClearFormats()
$ws.Range('E12').Value = '  +1.24%  '

$c = $ws.Range('D13')
$c.Value = "'1.627.31"
$c.ClearFormats()
$ws.Range('E13').Value = '  +1.12%  '

$ws.Range('E14').Value = '  +5.89%  '

$c = $ws.Range('D15')
$c.Value = "'3.90"
$c.ClearFormats()
$ws.Range('E15').Value = '  +5.04%  '

$c = $ws.Range('D16')
$c.Value = "'29.910.53"
$c.ClearFormats()
$ws.Range('E16').Value = '  +1.19%  '

$c = $ws.Range('D17')
$c.Value = "'8.79"
$c.ClearFormats()
$ws.Range('E17').Value = '  +15.91%  '

$c = $ws.Range('D18')
$c.Value = "'64.45"
$c.ClearFormats()
$ws.Range('E18').Value = '  +1.63%  '

$c = $ws.Range('D19')
$c.Value = "'243.59"
$c.ClearFormats()
$ws.Range('E19').Value = '  +1.34%  '

$c = $ws.Range('D20')
$c.Value = "'0.0₃0705"
$c.ClearFormats()

$ws.Range('E21').Value = '  -0.25%  '

$c = $ws.Range('D22')
$c.Value = "'4.12"
$c.ClearFormats()
$ws.Range('E22').Value = '  +3.49%  '

$c = $ws.Range('D23')
$c.Value = "'9.61"
$c.ClearFormats()
$ws.Range('E23').Value = '  +4.37%  '

$c = $ws.Range('D24')
$c.Value = "'2.13"
$c.ClearFormats()
$ws.Range('E24').Value = '  +2.28%  '

$c = $ws.Range('D25')
$c.Value = "'156.72"
$c.ClearFormats()
$ws.Range('E25').Value = '  +1.67%  '

$c = $ws.Range('D26')
$c.Value = "'15.64"
$c.ClearFormats()
$ws.Range('E26').Value = '  +2.58%  '

$ws.Range('E27').Value = '  +1.72%  '

$ws.Range('E28').Value = '  +2.95%  '

$c = $ws.Range('D29')
$c.Value = "'0.994"
$c.ClearFormats()
$ws.Range('E29').Value = '  -0.49%  '

$ws.Range('E30').Value = '  +3.60%  '

$c = $ws.Range('D31')
$c.Value = "'1.12"
$c.ClearFormats()
$ws.Range('E31').Value = '  +5.65%  '

$c = $ws.Range('D32')
$c.Value = "'3.34"
$c.ClearFormats()
$ws.Range('E32').Value = '  +3.60%  '

$ws.Range('E33').Value = '  +3.93%  '

$c = $ws.Range('D34')
$c.Value = "'1.426.03"
$c.ClearFormats()
$ws.Range('E34').Value = '  +1.35%  '

$ws.Range('E35').Value = '  +7.05%  '

$ws.Range('E36').Value = '  +0.00%  '

$ws.Range('E37').Value = '  +1.35%  '

$ws.Range('E38').Value = '  -0.71%  '

$ws.Range('E39').Value = '  +3.08%  '

$c = $ws.Range('D40')
$c.Value = "'0.555"
$c.ClearFormats()
$ws.Range('E40').Value = '  +3.17%  '

$ws.Range('E41').Value = '  +2.83%  '

$ws.Range('E42').Value = '  +0.12%  '

$c = $ws.Range('D43')
$c.Value = "'0.832"
$c.ClearFormats()
$ws.Range('E43').Value = '  +4.47%  '

$c = $ws.Range('D44')
$c.Value = "'54.32"
$c.ClearFormats()
$ws.Range('E44').Value = '  +0.95%  '

$c = $ws.Range('D45')
$c.Value = "'69.01"
$c.ClearFormats()
$ws.Range('E45').Value = '  +4.96%  '

$ws.Range('E46').Value = '  +19.38%  '

$c = $ws.Range('D47')
$c.Value = "'0.994"
$c.ClearFormats()
$ws.Range('E47').Value = '  -0.41%  '

$ws.Range('E48').Value = '  +2.68%  '

$c = $ws.Range('D49')
$c.Value = "'1.764.03"
$c.ClearFormats()
$ws.Range('E49').Value = '  +1.07%  '

$c = $ws.Range('D50')
$c.Value = "'88.06"
$c.ClearFormats()
$ws.Range('E50').Value = '  +1.72%  '

$c = $ws.Range('D51')
$c.Value = "'0.0₆0108"
$c.ClearFormats()
$ws.Range('E51').Value = '  +6.95%  '
